$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q4" right after "总计" (sheet 1),
#    built from a copy of the "2022-Q3" sheet so it inherits the same
#    sheet-level properties / column styles, then its cells are
#    overwritten with the 2022-Q4 fund-holding figures.
# ------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $firstSheet)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Helper: assign a value that must round-trip as *text* even though it
# looks numeric (fund codes / percentages stored as strings in the
# source data), using the apostrophe text-prefix like a user typing it.
function Set-TextValue($cell, [string]$value) {
    $cell.Formula = "'" + $value
}

# Row 2: 180028 银华永祥灵活配置混合
$q4.Cells.Item(2,1).Value = 0
Set-TextValue $q4.Cells.Item(2,2) "180028"
Set-TextValue $q4.Cells.Item(2,3) "银华永祥灵活配置混合"
Set-TextValue $q4.Cells.Item(2,4) "0.70"
Set-TextValue $q4.Cells.Item(2,5) "77.51"
Set-TextValue $q4.Cells.Item(2,6) "3.98"
Set-TextValue $q4.Cells.Item(2,7) "0.0279"
$q4.Cells.Item(2,8).Value = 8

# Row 3: 015694 瑞达策略优选混合A
$q4.Cells.Item(3,1).Value = 1
Set-TextValue $q4.Cells.Item(3,2) "015694"
Set-TextValue $q4.Cells.Item(3,3) "瑞达策略优选混合A"
Set-TextValue $q4.Cells.Item(3,4) "0.08"
Set-TextValue $q4.Cells.Item(3,5) "76.83"
Set-TextValue $q4.Cells.Item(3,6) "2.77"
Set-TextValue $q4.Cells.Item(3,7) "0.0022"
$q4.Cells.Item(3,8).Value = 6

# Row 4 (new row - copy column-A styling from the row above it first so
# the row-id cell keeps the same look as the rest of the sheet).
$q4.Cells.Item(3,1).Copy()
$q4.Cells.Item(4,1).PasteSpecial(-4122) # xlPasteFormats

$q4.Cells.Item(4,1).Value = 2
Set-TextValue $q4.Cells.Item(4,2) "015695"
Set-TextValue $q4.Cells.Item(4,3) "瑞达策略优选混合C"
Set-TextValue $q4.Cells.Item(4,4) "0.00"
Set-TextValue $q4.Cells.Item(4,5) "76.83"
Set-TextValue $q4.Cells.Item(4,6) "2.77"
$q4.Cells.Item(4,7).Value = 0
$q4.Cells.Item(4,8).Value = 6

# ------------------------------------------------------------------
# 2) Update the summary sheet "总计" (sheet 1): add a new first data
#    row for 2022-Q4 and push the existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy column-A styling down onto the freshly-needed row 7.
$summary.Cells.Item(6,1).Copy()
$summary.Cells.Item(7,1).PasteSpecial(-4122) # xlPasteFormats

# Existing rows 2-6 (2022-Q3 .. 2021-Q1) move down to rows 3-7. Values
# are written from the bottom up so nothing is overwritten before it is
# copied down to its new location.
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q1"
$summary.Cells.Item(7,3).Value = 4
$summary.Cells.Item(7,4).Value = 0.26

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q2"
$summary.Cells.Item(6,3).Value = 6
$summary.Cells.Item(6,4).Value = 0.59

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2022-Q1"
$summary.Cells.Item(5,3).Value = 2
$summary.Cells.Item(5,4).Value = 0.21

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q2"
$summary.Cells.Item(4,3).Value = 4
$summary.Cells.Item(4,4).Value = 0.17

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 2
$summary.Cells.Item(3,4).Value = 0.69

# New row 2: 2022-Q4
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.03

# Restore the originally-active sheet/selection (the copy/rename
# operations above shift the active tab onto the new sheet otherwise).
$wb.Worksheets.Item("2021-Q1").Activate()
[void]$wb.Worksheets.Item("2021-Q1").Range("A1").Select()

Write-Output "done"
